$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GREEK c)")

# Fill in the test results for "teste pasta 1", "teste pasta 2" and "teste pasta 3" rows
$ws.Range("B5").Value = 0.1
$ws.Range("C5").Value = 0.4
$ws.Range("D5").Value = 0.4

$ws.Range("B6").Value = 0.24
$ws.Range("C6").Value = 0.35
$ws.Range("D6").Value = 0.39

$ws.Range("B7").Value = 0.6
$ws.Range("C7").Value = 0.4
$ws.Range("D7").Value = 0.92500000000000004

# Add the same 3-colour colour-scale conditional formatting used on B2:D2
# to the newly filled rows, in sheet order, so priorities renumber the
# same way Excel does when each new rule is inserted with top priority.
$cf5 = $ws.Range("B5:D5").FormatConditions.AddColorScale(3)
$cf5.SetFirstPriority()

$cf6 = $ws.Range("B6:D6").FormatConditions.AddColorScale(3)
$cf6.SetFirstPriority()

$cf7 = $ws.Range("B7:D7").FormatConditions.AddColorScale(3)
$cf7.SetFirstPriority()

# Update the active selection on the sheet
$ws.Range("L5").Select()
